# edit.ps1 - apply the "minor change in ppt & dsa question" commit:
#  1. Fix the typo in the title slide: "Shell Short" -> "Shell Sort".
#  2. Tidy up the "Collision" label shape (drop its redundant / leftover
#     end-of-paragraph run formatting) on slide 2.
#  3. Add a new (empty, "Title and Content") slide 3 - content to be
#     filled in later, in keeping with the "less unnecessary content"
#     philosophy from the commit message.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Slide 1 (title slide): correct "Shell Short" -> "Shell Sort"
# ---------------------------------------------------------------------
$titleSlide = $p.Slides.Item(1)
$titleRange = $titleSlide.Shapes.Item(1).TextFrame.TextRange

$fullText = $titleRange.Text
$shellStart = $fullText.IndexOf("Shell ") + 1   # 1-based char index
$shellLen = "Shell ".Length
$shellPart = $titleRange.Characters($shellStart, $shellLen)
$shellPart.Text = "Shell "

$wrongStart = $fullText.IndexOf("Short") + 1    # 1-based char index
$wrongLen = "Short".Length
$wrongPart = $titleSlide.Shapes.Item(1).TextFrame.TextRange.Characters($wrongStart, $wrongLen)
$wrongPart.Text = "Sort"

# ---------------------------------------------------------------------
# 2) Slide 2: clean up the "Collision" shape inside the grouped arrow -
#    remove the stale duplicate end-paragraph formatting left over on
#    that text box.
# ---------------------------------------------------------------------
$introSlide = $p.Slides.Item(2)
$arrowGroup = $introSlide.Shapes.Item(3)
$collisionShape = $arrowGroup.GroupItems.Item(1)
$collisionRange = $collisionShape.TextFrame.TextRange
$collisionText = $collisionRange.Text
$collisionRange.Delete()
$collisionShape.TextFrame.TextRange.Text = $collisionText

# ---------------------------------------------------------------------
# 3) Add a new slide (slide 3) using the same "Title and Content" layout
#    as slide 2, left blank for now (content to follow later).
# ---------------------------------------------------------------------
$masterForNewSlide = $introSlide.Master
$titleAndContentLayout = $masterForNewSlide.CustomLayouts.Item(2)
$newSlide = $p.Slides.AddSlide($p.Slides.Count + 1, $titleAndContentLayout)
